{"js": "// Apply the within100.docx edit: update the date header and all 100\n// arithmetic-problem cells to their new values. Each entry is\n// [oldText, newText] in document order (date paragraph first, then the\n// 100 table-cell paragraphs, row-major, 20 rows x 5 columns). We walk\n// body.paragraphs (which enumerates paragraphs inside table cells too,\n// in document order), verify the existing text matches what we expect,\n// and replace it in place with insertText(..., Word.InsertLocation.replace)\n// so the run's formatting (fonts, size, etc.) is preserved.\n\nconst REPLACEMENTS = [\n  [\"2023-04-11 Tuesday\", \"2023-04-12 Wednesday\"], [\"28-11=\", \"47-44=\"], [\"6+8=\", \"22+63=\"], [\"42+21=\", \"13+29=\"],\n  [\"80-34=\", \"19+69=\"], [\"7+76=\", \"37+8=\"], [\"36+4=\", \"36-27=\"], [\"50+16=\", \"75-46=\"],\n  [\"78-3=\", \"15+1=\"], [\"63+17=\", \"61-13=\"], [\"89+5=\", \"37+46=\"], [\"40+7=\", \"27+33=\"],\n  [\"28-16=\", \"83-1=\"], [\"8+45=\", \"85+7=\"], [\"97-91=\", \"98-98=\"], [\"40+4=\", \"5+7=\"],\n  [\"81-68=\", \"73-38=\"], [\"87+9=\", \"86-34=\"], [\"98-82=\", \"20+6=\"], [\"55-40=\", \"38-24=\"],\n  [\"59-41=\", \"13+40=\"], [\"98-60=\", \"14+2=\"], [\"92-34=\", \"83-0=\"], [\"37-19=\", \"57+2=\"],\n  [\"56-46=\", \"96-38=\"], [\"87-20=\", \"19+76=\"], [\"48-3=\", \"70-12=\"], [\"94-9=\", \"36+61=\"],\n  [\"40+0=\", \"51+4=\"], [\"17-2=\", \"66+4=\"], [\"58-52=\", \"31+3=\"], [\"36-33=\", \"43-23=\"],\n  [\"78-62=\", \"28-17=\"], [\"87-15=\", \"54-24=\"], [\"30+69=\", \"56+19=\"], [\"9-0=\", \"69-41=\"],\n  [\"93-35=\", \"17+72=\"], [\"45+32=\", \"22+38=\"], [\"85-57=\", \"47+6=\"], [\"54-20=\", \"71-47=\"],\n  [\"41+47=\", \"8+83=\"], [\"62-4=\", \"87-47=\"], [\"69-28=\", \"97-22=\"], [\"85-5=\", \"82-5=\"],\n  [\"67-39=\", \"94-29=\"], [\"41+57=\", \"15+62=\"], [\"55+17=\", \"73-16=\"], [\"7+69=\", \"82-39=\"],\n  [\"99-67=\", \"86-41=\"], [\"90-46=\", \"49+7=\"], [\"17+6=\", \"41+55=\"], [\"90-14=\", \"58+13=\"],\n  [\"41+9=\", \"90-22=\"], [\"77+15=\", \"81-8=\"], [\"20+75=\", \"80+13=\"], [\"62-14=\", \"1+68=\"],\n  [\"41+33=\", \"85-3=\"], [\"70-58=\", \"11-4=\"], [\"91-33=\", \"97-42=\"], [\"21-9=\", \"81+6=\"],\n  [\"7+65=\", \"69+21=\"], [\"82-74=\", \"60+5=\"], [\"52+15=\", \"31+21=\"], [\"67-7=\", \"87-63=\"],\n  [\"75+12=\", \"49-17=\"], [\"17+61=\", \"55-26=\"], [\"20+73=\", \"63+9=\"], [\"26-5=\", \"55+43=\"],\n  [\"2+76=\", \"35+51=\"], [\"29+22=\", \"95-32=\"], [\"95-92=\", \"74-11=\"], [\"75+17=\", \"54+14=\"],\n  [\"70-0=\", \"48-4=\"], [\"92-41=\", \"64-15=\"], [\"39+47=\", \"57+25=\"], [\"95-62=\", \"36+58=\"],\n  [\"95-44=\", \"21-12=\"], [\"55-36=\", \"51-6=\"], [\"88-82=\", \"91-14=\"], [\"71+22=\", \"98-23=\"],\n  [\"63-28=\", \"96-76=\"], [\"97-91=\", \"51-16=\"], [\"57+38=\", \"11+77=\"], [\"5+76=\", \"42-17=\"],\n  [\"81-57=\", \"48+22=\"], [\"12-11=\", \"49+14=\"], [\"66-35=\", \"26+59=\"], [\"23+37=\", \"87-25=\"],\n  [\"64-14=\", \"22+71=\"], [\"17+64=\", \"77-1=\"], [\"33+55=\", \"42-14=\"], [\"90-57=\", \"44+47=\"],\n  [\"52+28=\", \"92-52=\"], [\"69+25=\", \"12+13=\"], [\"77+12=\", \"57-11=\"], [\"90-66=\", \"47-34=\"],\n  [\"97-75=\", \"13+82=\"], [\"49+18=\", \"45+54=\"], [\"52-5=\", \"17+19=\"], [\"49-26=\", \"50+23=\"],\n  [\"41-3=\", \"42-19=\"]\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== REPLACEMENTS.length) {\n  throw new Error(\n    `Expected ${REPLACEMENTS.length} paragraphs, found ${items.length}`\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const [oldText, newText] = REPLACEMENTS[i];\n  const para = items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      `Paragraph ${i}: expected \"${oldText}\" but found \"${para.text}\"`\n    );\n  }\n  if (oldText !== newText) {\n    para.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the within100.docx edit via Word COM interop:\n# - update the date header paragraph\n# - update each of the 100 arithmetic-problem table cells (20 rows x 5 cols)\n# to their new values, verifying the existing text first so we fail loudly\n# instead of silently mismatching if the document shape differs.\n\n$d = $word.ActiveDocument\n\n$dateOld = '2023-04-11 Tuesday'\n$dateNew = '2023-04-12 Wednesday'\n\n$dateRange = $d.Paragraphs.Item(1).Range\n$dateText = $dateRange.Text.TrimEnd([char]13, [char]7)\nif ($dateText -ne $dateOld) {\n    throw \"Date paragraph: expected `\"$dateOld`\" but found `\"$dateText`\"\"\n}\n$dateRange.Text = $dateNew\n\n# [oldText, newText] pairs for each table cell, row-major (20 rows x 5 cols),\n# matching the document's table exactly.\n$cellValues = @(\n    @(@('28-11=', '47-44='), @('6+8=', '22+63='), @('42+21=', '13+29='), @('80-34=', '19+69='), @('7+76=', '37+8=')),\n    @(@('36+4=', '36-27='), @('50+16=', '75-46='), @('78-3=', '15+1='), @('63+17=', '61-13='), @('89+5=', '37+46=')),\n    @(@('40+7=', '27+33='), @('28-16=', '83-1='), @('8+45=', '85+7='), @('97-91=', '98-98='), @('40+4=', '5+7=')),\n    @(@('81-68=', '73-38='), @('87+9=', '86-34='), @('98-82=', '20+6='), @('55-40=', '38-24='), @('59-41=', '13+40=')),\n    @(@('98-60=', '14+2='), @('92-34=', '83-0='), @('37-19=', '57+2='), @('56-46=', '96-38='), @('87-20=', '19+76=')),\n    @(@('48-3=', '70-12='), @('94-9=', '36+61='), @('40+0=', '51+4='), @('17-2=', '66+4='), @('58-52=', '31+3=')),\n    @(@('36-33=', '43-23='), @('78-62=', '28-17='), @('87-15=', '54-24='), @('30+69=', '56+19='), @('9-0=', '69-41=')),\n    @(@('93-35=', '17+72='), @('45+32=', '22+38='), @('85-57=', '47+6='), @('54-20=', '71-47='), @('41+47=', '8+83=')),\n    @(@('62-4=', '87-47='), @('69-28=', '97-22='), @('85-5=', '82-5='), @('67-39=', '94-29='), @('41+57=', '15+62=')),\n    @(@('55+17=', '73-16='), @('7+69=', '82-39='), @('99-67=', '86-41='), @('90-46=', '49+7='), @('17+6=', '41+55=')),\n    @(@('90-14=', '58+13='), @('41+9=', '90-22='), @('77+15=', '81-8='), @('20+75=', '80+13='), @('62-14=', '1+68=')),\n    @(@('41+33=', '85-3='), @('70-58=', '11-4='), @('91-33=', '97-42='), @('21-9=', '81+6='), @('7+65=', '69+21=')),\n    @(@('82-74=', '60+5='), @('52+15=', '31+21='), @('67-7=', '87-63='), @('75+12=', '49-17='), @('17+61=', '55-26=')),\n    @(@('20+73=', '63+9='), @('26-5=', '55+43='), @('2+76=', '35+51='), @('29+22=', '95-32='), @('95-92=', '74-11=')),\n    @(@('75+17=', '54+14='), @('70-0=', '48-4='), @('92-41=', '64-15='), @('39+47=', '57+25='), @('95-62=', '36+58=')),\n    @(@('95-44=', '21-12='), @('55-36=', '51-6='), @('88-82=', '91-14='), @('71+22=', '98-23='), @('63-28=', '96-76=')),\n    @(@('97-91=', '51-16='), @('57+38=', '11+77='), @('5+76=', '42-17='), @('81-57=', '48+22='), @('12-11=', '49+14=')),\n    @(@('66-35=', '26+59='), @('23+37=', '87-25='), @('64-14=', '22+71='), @('17+64=', '77-1='), @('33+55=', '42-14=')),\n    @(@('90-57=', '44+47='), @('52+28=', '92-52='), @('69+25=', '12+13='), @('77+12=', '57-11='), @('90-66=', '47-34=')),\n    @(@('97-75=', '13+82='), @('49+18=', '45+54='), @('52-5=', '17+19='), @('49-26=', '50+23='), @('41-3=', '42-19='))\n)\n\n$table = $d.Tables.Item(1)\n$numRows = $cellValues.Count\nfor ($r = 1; $r -le $numRows; $r++) {\n    $rowValues = $cellValues[$r - 1]\n    $numCols = $rowValues.Count\n    for ($c = 1; $c -le $numCols; $c++) {\n        $pair = $rowValues[$c - 1]\n        $oldText = $pair[0]\n        $newText = $pair[1]\n        $cell = $table.Cell($r, $c)\n        $cellRange = $cell.Range\n        $currentText = $cellRange.Text.TrimEnd([char]13, [char]7)\n        if ($currentText -ne $oldText) {\n            throw \"Cell ($r,$c): expected `\"$oldText`\" but found `\"$currentText`\"\"\n        }\n        $cellRange.Text = $newText\n    }\n}\n\n"}
